# Update gh-pages to output generated at 456a3b4
# Applies small numeric corrections to the "想去人数" (column F) values
# on the "展览" (sheet1) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 255
$ws1.Range("F23").Value = 10486
$ws1.Range("F27").Value = 2352
$ws1.Range("F29").Value = 2265
$ws1.Range("F34").Value = 2159
$ws1.Range("F46").Value = 1010
$ws1.Range("F47").Value = 1411
$ws1.Range("F48").Value = 71
$ws1.Range("F49").Value = 1107

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F27").Value = 10486
$ws4.Range("F29").Value = 2352
$ws4.Range("F31").Value = 2265
$ws4.Range("F36").Value = 2159
$ws4.Range("F48").Value = 1010
$ws4.Range("F49").Value = 1411
$ws4.Range("F50").Value = 71
$ws4.Range("F51").Value = 1107
